$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New score data for columns R:V (5 extra rating columns), rows 1-20
$data = @(
  @(5,5,5,4,5),
  @(3,3,3,3,3),
  @(1,1,2,1,2),
  @(2,3,2,2,2),
  @(4,4,4,4,5),
  @(3,2,3,2,3),
  @(1,1,1,1,2),
  @(2,1,3,1,2),
  @(5,5,4,5,5),
  @(3,3,4,3,3),
  @(1,1,1,1,1),
  @(3,3,3,2,3),
  @(4,5,5,4,4),
  @(2,2,3,2,2),
  @(1,1,2,1,2),
  @(2,3,3,2,3),
  @(3,5,5,4,4),
  @(2,2,3,2,3),
  @(1,1,2,1,1),
  @(2,3,3,2,2)
)

# Copy the formatting already used on column Q (style index matches the
# rest of the data columns) onto the new R:V range before filling values.
$ws.Range("Q1:Q20").Copy()
$ws.Range("R1:V20").PasteSpecial(-4122)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = 18 + $j   # R=18, S=19, T=20, U=21, V=22
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}

# Match the final cell selection left by the author after entering the data
$ws.Range("V21").Select()
